$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 16, pushing the existing data (old rows 16-98)
# down to rows 18-100.
$ws.Range("A16:A17").EntireRow.Insert()

# --- New row 16 -------------------------------------------------------
$ws.Cells.Item(16, 1).Value = 5
$ws.Cells.Item(16, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(16, 3).Value = "Maule"
$ws.Cells.Item(16, 4).Value = 45114
$ws.Cells.Item(16, 5).Value = 7
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100104
$ws.Cells.Item(16, 8).Value = "Frutos de pepita"
$ws.Cells.Item(16, 9).Value = 100104003
$ws.Cells.Item(16, 10).Value = "Membrillo"
$ws.Cells.Item(16, 11).Value = "Champion"
$ws.Cells.Item(16, 12).Value = "Especial"
$ws.Cells.Item(16, 13).Value = 200
$ws.Cells.Item(16, 14).Value = 12000
$ws.Cells.Item(16, 15).Value = 12000
$ws.Cells.Item(16, 16).Value = 12000
$ws.Cells.Item(16, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(16, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(16, 19).Value = 667
$ws.Cells.Item(16, 20).Value = 18

# --- New row 17 -------------------------------------------------------
$ws.Cells.Item(17, 1).Value = 5
$ws.Cells.Item(17, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(17, 3).Value = "Maule"
$ws.Cells.Item(17, 4).Value = 45114
$ws.Cells.Item(17, 5).Value = 7
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100104
$ws.Cells.Item(17, 8).Value = "Frutos de pepita"
$ws.Cells.Item(17, 9).Value = 100104003
$ws.Cells.Item(17, 10).Value = "Membrillo"
$ws.Cells.Item(17, 11).Value = "Champion"
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 200
$ws.Cells.Item(17, 14).Value = 10000
$ws.Cells.Item(17, 15).Value = 10000
$ws.Cells.Item(17, 16).Value = 10000
$ws.Cells.Item(17, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(17, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(17, 19).Value = 556
$ws.Cells.Item(17, 20).Value = 18
